$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.048.40"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "1.886.92"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7377"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.24"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3170"
$ws.Range("E8").Value = "  +2.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07187"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.78"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08332"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7581"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").Value = "1.898.12"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.398"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.12"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").Value = "30.052.14"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "250.36"
$ws.Range("E18").Value = "  +4.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.57"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.140.85"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.896"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1560"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.12"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.70"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.049"
$ws.Range("E29").Value = "  +2.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.479"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.567"
$ws.Range("E31").Value = "  +3.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.533"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.206"
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05341"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.251"
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7701"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9986"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.720"
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.759"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4577"
$ws.Range("E41").Value = "  +3.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.046"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").Value = "1.090.53"
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.35"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8728"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.80"
$ws.Range("E46").Value = "  +2.97%  "
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.858"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.587"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.558"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "2.045.23"
$ws.Range("E51").Value = "  +0.53%  "
